$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Model Accuracy (-0.8, 0.8, 0.8)"
#   - Add new columns C:G (Market threshold, Market min, Market max, Recall,
#     Precision) with header row + 5 data rows.
#   - Update the existing Accuracy (%) values in column B.
# ---------------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item(1)

# New header cells (C1:G1) - copy formatting from the existing B1 header cell
$wsAcc.Range("C1").Value = "Market threshold"
$wsAcc.Range("D1").Value = "Market min"
$wsAcc.Range("E1").Value = "Market max"
$wsAcc.Range("F1").Value = "Recall"
$wsAcc.Range("G1").Value = "Precision"

$wsAcc.Range("B1").Copy()
$wsAcc.Range("C1:G1").PasteSpecial(-4122)

# Row 2 - TOTALENERGIES SE
$wsAcc.Range("B2").Value = 63.20293398533008
$wsAcc.Range("C2").Value = 0.05450546436368681
$wsAcc.Range("D2").Value = -15.55441
$wsAcc.Range("E2").Value = 15.06418
$wsAcc.Range("F2").Value = 0
$wsAcc.Range("G2").Value = 0

# Row 3 - FMC CORP
$wsAcc.Range("B3").Value = 38.14180929095355
$wsAcc.Range("C3").Value = 0.009583939973006913
$wsAcc.Range("D3").Value = -19.35264
$wsAcc.Range("E3").Value = 13.70093
$wsAcc.Range("F3").Value = 2.144772117962467
$wsAcc.Range("G3").Value = 26.66666666666667

# Row 4 - BP PLC
$wsAcc.Range("B4").Value = 92.66503667481662
$wsAcc.Range("C4").Value = 0.04158117063764853
$wsAcc.Range("D4").Value = -18.75314
$wsAcc.Range("E4").Value = 23.33066
$wsAcc.Range("F4").Value = 0
$wsAcc.Range("G4").Value = 0

# Row 5 - STORA ENSO
$wsAcc.Range("B5").Value = 82.09046454767727
$wsAcc.Range("C5").Value = 0.02983403801513819
$wsAcc.Range("D5").Value = -12.78028
$wsAcc.Range("E5").Value = 12.42348
$wsAcc.Range("F5").Value = 0
$wsAcc.Range("G5").Value = 0

# Row 6 - BHP GROUP
$wsAcc.Range("B6").Value = 95.59902200488997
$wsAcc.Range("C6").Value = 0.08368817696170747
$wsAcc.Range("D6").Value = -16.47904
$wsAcc.Range("E6").Value = 14.94325
$wsAcc.Range("F6").Value = 0
$wsAcc.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Confusion Matrix TOTALENERGIES SE (-0.8, 0.8, 0.8)"
#   - Row 3 (Predicted Neutral) values change.
# ---------------------------------------------------------------------------
$wsCm1 = $wb.Worksheets.Item(2)
$wsCm1.Range("B3").Value = 9
$wsCm1.Range("C3").Value = 1033
$wsCm1.Range("D3").Value = 9

# ---------------------------------------------------------------------------
# Sheet 3: "Confusion Matrix FMC CORP (-0.8, 0.8, 0.8)"
#   - Rows 2, 3 and 4 values change.
# ---------------------------------------------------------------------------
$wsCm2 = $wb.Worksheets.Item(3)
$wsCm2.Range("B2").Value = 8
$wsCm2.Range("C2").Value = 17
$wsCm2.Range("D2").Value = 5

$wsCm2.Range("B3").Value = 339
$wsCm2.Range("C3").Value = 589
$wsCm2.Range("D3").Value = 324

$wsCm2.Range("B4").Value = 26
$wsCm2.Range("C4").Value = 46
$wsCm2.Range("D4").Value = 27

# ---------------------------------------------------------------------------
# Sheet 4: "Confusion Matrix BP PLC (-0.8, 0.8, 0.8)"
#   - Row 3 (Predicted Neutral) values change.
# ---------------------------------------------------------------------------
$wsCm3 = $wb.Worksheets.Item(4)
$wsCm3.Range("B3").Value = 40
$wsCm3.Range("C3").Value = 1516
$wsCm3.Range("D3").Value = 42

# ---------------------------------------------------------------------------
# Sheet 5: "Confusion Matrix STORA ENSO (-0.8, 0.8, 0.8)"
#   - Row 3 (Predicted Neutral) values change.
# ---------------------------------------------------------------------------
$wsCm4 = $wb.Worksheets.Item(5)
$wsCm4.Range("B3").Value = 110
$wsCm4.Range("C3").Value = 1343
$wsCm4.Range("D3").Value = 107

# ---------------------------------------------------------------------------
# Sheet 6: "Confusion Matrix BHP GROUP (-0.8, 0.8, 0.8)"
#   - Row 3 (Predicted Neutral) values change.
# ---------------------------------------------------------------------------
$wsCm5 = $wb.Worksheets.Item(6)
$wsCm5.Range("B3").Value = 4
$wsCm5.Range("C3").Value = 1564
$wsCm5.Range("D3").Value = 3

# Sheet 7 "Signal Market Correlation (-0.8, 0.8, 0.8)" is unchanged.
